# "Put base registration rate as separate sheet"
#
# This adds a new variable ("Base registration rate" / RSHORTTI family row)
# to the FTT-Tr sheet (as a new row 8) and registers it in the Time_Horizons
# sheet (as a new row 10), shifting the subsequent rows down on both sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. FTT-Tr sheet ("Inputs" table): insert the new variable as row 8
# ---------------------------------------------------------------------
$wsTr = $wb.Worksheets.Item(2)

$wsTr.Rows("8:8").Insert()

$wsTr.Cells.Item(8, 1).Value = "Base registration rate"
$wsTr.Cells.Item(8, 2).Value = 1
$wsTr.Cells.Item(8, 3).Value = -99
$wsTr.Cells.Item(8, 4).Value = "FTT-Tr relative registration tax or subsidy"
$wsTr.Cells.Item(8, 5).Value = "VTTI"
$wsTr.Cells.Item(8, 6).Value = "TIME"
$wsTr.Cells.Item(8, 7).Value = "RSHORTTI"
$wsTr.Cells.Item(8, 8).Value = 0
$wsTr.Cells.Item(8, 9).Value = "All"

# Widen column D so the longer description text fits.
$wsTr.Columns.Item(4).ColumnWidth = 45.86

# ---------------------------------------------------------------------
# 2. Time_Horizons sheet: register the new variable as row 10
# ---------------------------------------------------------------------
$wsTh = $wb.Worksheets.Item(6)

$wsTh.Rows("10:10").Insert()

$wsTh.Cells.Item(10, 1).Value = "Base registration rate"
$wsTh.Cells.Item(10, 2).Value = "tl_2001"

# ---------------------------------------------------------------------
# 3. Restore/update the view state (selections + active sheet/tab)
#    FTT-P's own selection (C16) is left untouched - only the active-tab
#    bookkeeping changes because Time_Horizons becomes the active sheet.
# ---------------------------------------------------------------------
$wsTr.Range("A3").Select()

$wsTh.Range("A11").Select()
$wsTh.Activate()
